# Purple Tentacle: Check for close actors (script-43)
# A new Variable (#74 "Weird Ed: Should storm lab") is inserted into the
# "Variables" sheet, right between variable #73 ("Door Bell Triggers Ed",
# row 53) and the former variable #75 ("Record on Victrola", old row 54).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Insert a new blank row at row 54, shifting everything below down by one.
$ws.Rows.Item(54).Insert()

# Fill in the new variable's data.
$ws.Cells.Item(54, 1).Value = 74
$ws.Cells.Item(54, 3).Value = "Weird Ed: Should storm lab"
$ws.Cells.Item(54, 5).Value = "1 = Yes"

# Leave the cursor where the author left it: Variables sheet, scrolled so
# row 40 is at the top, with E55 selected and active.
$ws.Application.Goto($ws.Range("A40"))
$ws.Range("E55").Select()

$scripts = $wb.Worksheets.Item("Scripts")
$scripts.Application.Goto($scripts.Range("A91"))
$scripts.Range("D93").Select()

$objects = $wb.Worksheets.Item("Objects")
$objects.Range("D23").Select()

$ws.Activate()
